# Auto-update files from Google Sheet trigger
#
# 1. product_listing!M6 gets the full SOCAT description (the "Description"
#    column was empty for the SOCAT row before this edit).
# 2. web_formatted!V6 ("card-detail.content") is repointed from
#    product_listing!L6 (Highlights) to product_listing!M6 (Description) so
#    the detail panel shows the new long-form text instead of the old
#    one-liner.
# 3. web_formatted!Q60:Q124 ("card-summary" tail rows, including the rows
#    that had degraded into #REF! shared formulas after earlier column
#    edits) are cleared out entirely.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("product_listing")
$ws3 = $wb.Worksheets.Item("web_formatted")

# Match the formatting already used by the surrounding cells in column M
# (e.g. M5) before putting the new text in M6.
$ws1.Range("M5").Copy()
$ws1.Range("M6").PasteSpecial(-4122)

$ws1.Range("M6").Value = "The Surface Ocean CO2 Atlas features surface fCO2 measurements from both the open ocean and the coastal ocean, predominantly sourced from research vessels, ships of opportunity, and autonomous platforms including fixed moorings and uncrewed surface vehicles (USVs) (Bakker et al., 2016). It represents the most extensive collection of observational ocean CO2 data for the global surface ocean. Since 2013, SOCAT has been updated annually. Dataset flags indicate the estimated uncertainty and completeness of metadata in SOCAT synthesis products. The SOCAT gridded product contains fCO2 values with an estimated uncertainty of less than 5 µatm. "

$ws3.Range("V6").Formula = "=product_listing!M6"

$ws3.Range("Q60:Q124").ClearContents()
